# Update countries & provincias Spain
# Applies the data refresh captured in the diff:
#  - Re-sorted / corrected country names for several rows (block of countries
#    around "San Vicente y las Granadinas", "Sierra Leona", "San Bartolome",
#    "Anguila", "Islas Virgenes Britanicas", "Papua Nueva Guinea",
#    "Islas Malvinas", "Bonaire, San Eustaquio y Saba", "Burundi", "Sudan del Sur")
#  - Updated case statistics for a handful of rows
#  - Updated the "Datos actualizados..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Title / timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 13:52"

# --- Country name corrections (column A) ---
$ws.Range("A191").Value = "San Vicente y las Granadinas"
$ws.Range("A192").Value = "Islas Turcas y Caicos"
$ws.Range("A193").Value = "Somalia"
$ws.Range("A195").Value = "Sierra Leona"
$ws.Range("A196").Value = "Santa Sede"
$ws.Range("A197").Value = "Belice"
$ws.Range("A198").Value = "Cabo Verde"
$ws.Range("A200").Value = "San Bartolome"
$ws.Range("A201").Value = "Botsuana"
$ws.Range("A207").Value = "Anguila"
$ws.Range("A209").Value = "Burundi"
$ws.Range("A210").Value = "Papua Nueva Guinea"
$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A213").Value = "Sudan del Sur"

# --- Statistic updates ---
# Polonia (row 31)
$ws.Range("F31").Value = 160

# Finlandia (row 43)
$ws.Range("E43").Value = 2147
$ws.Range("F43").Value = 82
$ws.Range("G43").Value = 6
$ws.Range("H43").Value = 40

# Kazajistan (row 76)
$ws.Range("D76").Value = 54
$ws.Range("E76").Value = 648

# Row 191 (San Vicente y las Granadinas)
$ws.Range("D191").Value = 1
$ws.Range("H191").Value = 0

# Row 192 (Islas Turcas y Caicos)
$ws.Range("D192").Value = 0
$ws.Range("H192").Value = 1

# Row 195 (Sierra Leona)
$ws.Range("C195").Value = 1

# Row 196 (Santa Sede)
$ws.Range("E196").Value = 7
$ws.Range("F196").Value = 0
$ws.Range("H196").Value = 0

# Row 197 (Belice)
$ws.Range("D197").Value = 0
$ws.Range("E197").Value = 6
$ws.Range("F197").Value = 1

# Row 198 (Cabo Verde)
$ws.Range("B198").Value = 7
$ws.Range("D198").Value = 1
$ws.Range("E198").Value = 5
$ws.Range("H198").Value = 1

# Row 200 (San Bartolome)
$ws.Range("D200").Value = 1
$ws.Range("H200").Value = 0

# Row 201 (Botsuana)
$ws.Range("D201").Value = 0
$ws.Range("H201").Value = 1
